$wb = $excel.ActiveWorkbook

# --- Update metadata on the "Metadata" sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "2.0.0"
$meta.Range("B8").Value = "2024-06-04T14:59:10+02:00"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Add a new "Include from FSIII 12" sheet at the end, mirroring the
#     layout of the previous "Include from FSIII N" sheets ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "Include from FSIII 12"

# Copy the layout/styles/column widths from the previous FSIII sheet.
$lastSheet.Range("A1:C4").Copy($newSheet.Range("A1"))

# Fill in the new concept's unique identifier and clear the leftover
# copied value in column C for rows 3/4 (those rows only use A/B).
$newSheet.Range("C2").Value = "aec684bd-c2ea-4ff0-8eb7-6d2cf67fb863"
$newSheet.Range("C3").Value = [System.Type]::Missing
$newSheet.Range("C4").Value = [System.Type]::Missing

$newSheet.Columns.Item(1).ColumnWidth = 30.703125
$newSheet.Columns.Item(2).ColumnWidth = 50.703125
